$wb = $excel.ActiveWorkbook

# Status text used once a locale has been handed back (was "Ready for handoff").
# This string is shared across the Overview roll-up sheet as well as each
# locale's detail sheet, so every cell currently showing "Ready for handoff"
# needs to move to the new text together.
$statusHandedBack = "Handed back: in sync with en-US"

# ---- Overview sheet (roll-up Status columns for each locale) -----------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("B2").Value = $statusHandedBack
$ws.Range("C2").Value = $statusHandedBack
$ws.Range("B3").Value = $statusHandedBack
$ws.Range("C3").Value = $statusHandedBack

# ---- zh-cn sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

# Status column -> handed back
$ws.Range("B2").Value = $statusHandedBack
$ws.Range("B3").Value = $statusHandedBack

# Latest Target File (E) / Latest Handback File (F) for the two content rows,
# mirroring the existing Latest Handoff File (A) / Latest Handoff File (C) links.
$ws.Range("E2").Value = "3d6d923a-4f6f-4169-992f-ccc384019ff3.md"
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/84c4155b4f44cb7535b982284fb95a6915a83b8c/e2e/3d6d923a-4f6f-4169-992f-ccc384019ff3.md", "", "", "3d6d923a-4f6f-4169-992f-ccc384019ff3.md") | Out-Null

$ws.Range("F2").Value = "3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.zh-cn.xlf"
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8936cbf64460a8d88a1cbecad0156059d56b583f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.zh-cn.xlf", "", "", "3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.zh-cn.xlf") | Out-Null

$ws.Range("E3").Value = "3d6d923a-4f6f-4169-992f-ccc384019ff3.md"
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/84c4155b4f44cb7535b982284fb95a6915a83b8c/e2e/3d6d923a-4f6f-4169-992f-ccc384019ff3.md", "", "", "3d6d923a-4f6f-4169-992f-ccc384019ff3.md") | Out-Null

$ws.Range("F3").Value = "3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.zh-cn.xlf"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8936cbf64460a8d88a1cbecad0156059d56b583f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.zh-cn.xlf", "", "", "3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.zh-cn.xlf") | Out-Null

# Latest Handback DateTime (G) for the two content rows.
$ws.Range("G2").Value = "2016-03-09 14:29:36"
$ws.Range("G3").Value = "2016-03-09 14:29:36"

# ---- de-de sheet --------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

# Status column -> handed back
$ws.Range("B2").Value = $statusHandedBack
$ws.Range("B3").Value = $statusHandedBack

# Latest Target File (E) / Latest Handback File (F) for the two content rows.
$ws.Range("E2").Value = "3d6d923a-4f6f-4169-992f-ccc384019ff3.md"
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/84c4155b4f44cb7535b982284fb95a6915a83b8c/e2e/3d6d923a-4f6f-4169-992f-ccc384019ff3.md", "", "", "3d6d923a-4f6f-4169-992f-ccc384019ff3.md") | Out-Null

$ws.Range("F2").Value = "3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.de-de.xlf"
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9294d13ee7f6855343eb677d18ff7b6b8dc09f84/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.de-de.xlf", "", "", "3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.de-de.xlf") | Out-Null

$ws.Range("E3").Value = "3d6d923a-4f6f-4169-992f-ccc384019ff3.md"
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/84c4155b4f44cb7535b982284fb95a6915a83b8c/e2e/3d6d923a-4f6f-4169-992f-ccc384019ff3.md", "", "", "3d6d923a-4f6f-4169-992f-ccc384019ff3.md") | Out-Null

$ws.Range("F3").Value = "3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.de-de.xlf"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9294d13ee7f6855343eb677d18ff7b6b8dc09f84/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.de-de.xlf", "", "", "3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.de-de.xlf") | Out-Null

# Latest Handback DateTime (G) for the two content rows.
$ws.Range("G2").Value = "2016-03-09 14:29:41"
$ws.Range("G3").Value = "2016-03-09 14:29:41"
